$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Sheet1: recomputed "RF concerntration" numeric values in column E (rows 2-13)
$ws1.Range("E2").Value = 187.4662394188178
$ws1.Range("E3").Value = 44.78663001715602
$ws1.Range("E4").Value = 444.4066809698979
$ws1.Range("E5").Value = 223.7151910842536
$ws1.Range("E6").Value = 30.4925135402917
$ws1.Range("E7").Value = 6.049374470398623
$ws1.Range("E8").Value = 91.25213241309717
$ws1.Range("E9").Value = 260.8309210117467
$ws1.Range("E10").Value = 191.5763917159027
$ws1.Range("E11").Value = 40.18389530692563
$ws1.Range("E12").Value = 946.7144485283916
$ws1.Range("E13").Value = 33.39695732022464

# Sheet2: "RF concentration" (E) values must stay text; force text format before writing
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "187.4662394188178"
$ws2.Range("B4").Value = "Best Match: C(C(C(C(F)(F)Cl)(F)F)(F)F)(C(C(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(F)F)(F)F with Similarity: 0.0740740740740741"
$ws2.Range("E6").NumberFormat = "@"
$ws2.Range("E6").Value = "44.78663001715602"
$ws2.Range("B8").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.0612244897959184"
$ws2.Range("E10").NumberFormat = "@"
$ws2.Range("E10").Value = "444.40668096989793"
$ws2.Range("B12").Value = "Best Match: C(C(C(C(F)(F)Cl)(F)F)(F)F)(C(C(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(F)F)(F)F with Similarity: 0.0588235294117647"
$ws2.Range("E14").NumberFormat = "@"
$ws2.Range("E14").Value = "223.7151910842536"
$ws2.Range("B16").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.0512820512820513"
$ws2.Range("E18").NumberFormat = "@"
$ws2.Range("E18").Value = "30.492513540291704"
$ws2.Range("B20").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.0384615384615385"
$ws2.Range("E22").NumberFormat = "@"
$ws2.Range("E22").Value = "6.049374470398623"
$ws2.Range("B24").Value = "Best Match: C(C(F)(F)F)(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F with Similarity: 0.028169014084507"
$ws2.Range("E26").NumberFormat = "@"
$ws2.Range("E26").Value = "91.25213241309717"
$ws2.Range("B28").Value = "Best Match: C(C(C(C(F)(F)Cl)(F)F)(F)F)(C(C(OC(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(F)F)(F)F with Similarity: 0.0510204081632653"
$ws2.Range("E30").NumberFormat = "@"
$ws2.Range("E30").Value = "260.8309210117467"
$ws2.Range("B32").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.0444444444444444"
$ws2.Range("E34").NumberFormat = "@"
$ws2.Range("E34").Value = "191.57639171590267"
$ws2.Range("B36").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.025974025974026"
$ws2.Range("E38").NumberFormat = "@"
$ws2.Range("E38").Value = "40.18389530692563"
$ws2.Range("B40").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.0222222222222222"
$ws2.Range("E42").NumberFormat = "@"
$ws2.Range("E42").Value = "946.7144485283916"
$ws2.Range("E46").NumberFormat = "@"
$ws2.Range("E46").Value = "33.39695732022464"
$ws2.Range("B48").Value = "Best Match: C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(C(C(F)(F)F)(F)F)(F)F with Similarity: 0.025974025974026"
